# Apply updated Betfair back/lay odds values to Sheet1
# (values refreshed for 2026-01-05 games)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 1.39
$ws.Range("I2").Value = 1.4
$ws.Range("J2").Value = 4.7
$ws.Range("K2").Value = 4.8
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 4.5
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 1.8
$ws.Range("Q2").Value = 2.22
$ws.Range("R2").Value = 1.22
$ws.Range("S2").Value = 5.4
$ws.Range("T2").Value = 1.96
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 3.5
$ws.Range("W2").Value = 1.07
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 4.9
$ws.Range("Z2").Value = 5.7
$ws.Range("AA2").Value = 14.5
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 6.4
$ws.Range("AD2").Value = 7.8
$ws.Range("AE2").Value = 22
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 21
$ws.Range("AH2").Value = 25
$ws.Range("AI2").Value = 70
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 140
$ws.Range("AL2").Value = 160
$ws.Range("AM2").Value = 370
$ws.Range("AN2").Value = 410
$ws.Range("AO2").Value = 25

# Row 3
$ws.Range("F3").Value = 3.7
$ws.Range("H3").Value = 2.06
$ws.Range("I3").Value = 2.18
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 3.95
$ws.Range("L3").Value = 1.46
$ws.Range("Q3").Value = 2.06
$ws.Range("R3").Value = 1.33
$ws.Range("S3").Value = 3.7
$ws.Range("T3").Value = 1.9
$ws.Range("U3").Value = 2.04
$ws.Range("V3").Value = 1.85
$ws.Range("W3").Value = 1.32
$ws.Range("X3").Value = 16.5
$ws.Range("Z3").Value = 16
$ws.Range("AA3").Value = 32
$ws.Range("AB3").Value = 14
$ws.Range("AC3").Value = 8.6
$ws.Range("AD3").Value = 990
$ws.Range("AE3").Value = 28
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 42
$ws.Range("AL3").Value = 70
$ws.Range("AN3").Value = 65
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 2.64
$ws.Range("G4").Value = 2.92
$ws.Range("H4").Value = 2.72
$ws.Range("I4").Value = 2.96
$ws.Range("K4").Value = 3.65
$ws.Range("L4").Value = 1.42
$ws.Range("N4").Value = 3.95
$ws.Range("P4").Value = 1.98
$ws.Range("Q4").Value = 1.91
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.18
$ws.Range("V4").Value = 1.52
$ws.Range("W4").Value = 1.52
$ws.Range("X4").Value = 990
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AC4").Value = 8.6
$ws.Range("AD4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AN4").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.45
$ws.Range("G5").Value = 1.51
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4.5
$ws.Range("L5").Value = 1.48
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 2.94
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 1.67
$ws.Range("Q5").Value = 2.34
$ws.Range("R5").Value = 1.23
$ws.Range("S5").Value = 4.6
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = 1.55
$ws.Range("V5").Value = 1.08
$ws.Range("W5").Value = 2.92
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 29
$ws.Range("Z5").Value = 140
$ws.Range("AB5").Value = 6.2
$ws.Range("AC5").Value = 11
$ws.Range("AD5").Value = 65
$ws.Range("AF5").Value = 7.6
$ws.Range("AH5").Value = 44
$ws.Range("AJ5").Value = 13.5
$ws.Range("AK5").Value = 22
$ws.Range("AL5").Value = 1000
$ws.Range("AN5").Value = 12.5

# Row 6
$ws.Range("F6").Value = 1.9
$ws.Range("G6").Value = 1.95
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 4.1
$ws.Range("N6").Value = 4.3
$ws.Range("P6").Value = 2.16
$ws.Range("R6").Value = 1.43
$ws.Range("S6").Value = 3.15
$ws.Range("U6").Value = 2.14
$ws.Range("V6").Value = 1.29
$ws.Range("W6").Value = 2.06
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 17.5
$ws.Range("Z6").Value = 34
$ws.Range("AA6").Value = 300
$ws.Range("AB6").Value = 10
$ws.Range("AC6").Value = 8.6
$ws.Range("AD6").Value = 17.5
$ws.Range("AE6").Value = 55
$ws.Range("AF6").Value = 12.5
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 18
$ws.Range("AI6").Value = 60
$ws.Range("AJ6").Value = 22
$ws.Range("AK6").Value = 19.5
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 200
$ws.Range("AN6").Value = 12.5
$ws.Range("AO6").Value = 48

# Row 7
$ws.Range("J7").Value = 5.2
$ws.Range("R7").Value = 1.74
$ws.Range("U7").Value = 2.18
$ws.Range("AF7").Value = 11.5
$ws.Range("AK7").Value = 14.5
$ws.Range("AN7").Value = 5.2

# Row 8
$ws.Range("F8").Value = 1.79
$ws.Range("G8").Value = 1.85
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 7.4
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 3.5
$ws.Range("L8").Value = 1.67
$ws.Range("M8").Value = 1.15
$ws.Range("N8").Value = 2.36
$ws.Range("O8").Value = 1.64
$ws.Range("P8").Value = 1.44
$ws.Range("Q8").Value = 2.92
$ws.Range("R8").Value = 1.15
$ws.Range("S8").Value = 6.4
$ws.Range("T8").Value = 2.56
$ws.Range("U8").Value = 1.55
$ws.Range("V8").Value = 1.16
$ws.Range("W8").Value = 2.16
$ws.Range("X8").Value = 8
$ws.Range("Y8").Value = 15
$ws.Range("AB8").Value = 5.5
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AD8").Value = 85
$ws.Range("AF8").Value = 9
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 85
$ws.Range("AJ8").Value = 21
$ws.Range("AK8").Value = 65
$ws.Range("AL8").Value = 1000
$ws.Range("AN8").Value = 600

# Row 9
$ws.Range("G9").Value = 1.86
$ws.Range("H9").Value = 4.9
$ws.Range("I9").Value = 5.7
$ws.Range("J9").Value = 3.85
$ws.Range("K9").Value = 4.1
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 3.35
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 1.79
$ws.Range("Q9").Value = 2.08
$ws.Range("T9").Value = 1.94
$ws.Range("W9").Value = 2.16
$ws.Range("AC9").Value = 9.800000000000001
$ws.Range("AF9").Value = 18
$ws.Range("AJ9").Value = 900

# Row 10
$ws.Range("F10").Value = 1.66
$ws.Range("G10").Value = 1.74
$ws.Range("H10").Value = 6.2
$ws.Range("I10").Value = 7.2
$ws.Range("K10").Value = 4.2
$ws.Range("L10").Value = 1.43
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 3.3
$ws.Range("O10").Value = 1.39
$ws.Range("P10").Value = 1.78
$ws.Range("Q10").Value = 2.14
$ws.Range("R10").Value = 1.29
$ws.Range("S10").Value = 4
$ws.Range("T10").Value = 2.06
$ws.Range("U10").Value = 1.75
$ws.Range("W10").Value = 2.32
$ws.Range("X10").Value = 27
$ws.Range("Y10").Value = 980
$ws.Range("AB10").Value = 14.5
$ws.Range("AF10").Value = 22
$ws.Range("AJ10").Value = 180

# Row 11
$ws.Range("F11").Value = 2.16
$ws.Range("G11").Value = 2.3
$ws.Range("L11").Value = 1.41
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 3.6
$ws.Range("O11").Value = 1.31
$ws.Range("P11").Value = 1.92
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 1.35
$ws.Range("T11").Value = 1.74
$ws.Range("U11").Value = 2.08
$ws.Range("Y11").Value = 14.5
$ws.Range("AA11").Value = 900
$ws.Range("AB11").Value = 10.5
$ws.Range("AH11").Value = 18.5
$ws.Range("AI11").Value = 330
$ws.Range("AL11").Value = 100
$ws.Range("AM11").Value = 580
$ws.Range("AN11").Value = 18
$ws.Range("AO11").Value = 600

# Row 12
$ws.Range("F12").Value = 2.44
$ws.Range("G12").Value = 2.56
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.4
$ws.Range("K12").Value = 3.35
$ws.Range("L12").Value = 1.49
$ws.Range("M12").Value = 1.09
$ws.Range("N12").Value = 3.3
$ws.Range("O12").Value = 1.4
$ws.Range("P12").Value = 1.75
$ws.Range("Q12").Value = 2.24
$ws.Range("R12").Value = 1.28
$ws.Range("S12").Value = 4.2
$ws.Range("T12").Value = 1.88
$ws.Range("U12").Value = 2.02
$ws.Range("V12").Value = 1.41
$ws.Range("W12").Value = 1.64
$ws.Range("X12").Value = 11.5
$ws.Range("Y12").Value = 11.5
$ws.Range("Z12").Value = 22
$ws.Range("AA12").Value = 60
$ws.Range("AB12").Value = 9.199999999999999
$ws.Range("AC12").Value = 7.4
$ws.Range("AD12").Value = 14
$ws.Range("AE12").Value = 44
$ws.Range("AH12").Value = 19
$ws.Range("AI12").Value = 70
$ws.Range("AJ12").Value = 36
$ws.Range("AK12").Value = 30
$ws.Range("AL12").Value = 46
$ws.Range("AM12").Value = 120
$ws.Range("AN12").Value = 28
$ws.Range("AO12").Value = 55

# Row 13
$ws.Range("F13").Value = 2
$ws.Range("I13").Value = 4.6
$ws.Range("J13").Value = 3.45
$ws.Range("O13").Value = 1.37
$ws.Range("S13").Value = 4
$ws.Range("T13").Value = 1.91
$ws.Range("AB13").Value = 9
$ws.Range("AF13").Value = 27

Write-Host "Applied 275 cell value updates"
